$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the typo in TC_ERP_Login_002's test case description (row 11, col D)
# ---------------------------------------------------------------------------
$ws.Range("D11").Value2 = "Enter valid user name and one invalid password"

# ---------------------------------------------------------------------------
# 2. Correct the execution date for TC_ERP_Login_002 (row 11, col M)
#    45114 (07-Jul-2023) -> 45117 (10-Jul-2023)
# ---------------------------------------------------------------------------
$ws.Range("M11").Value2 = 45117

# ---------------------------------------------------------------------------
# 3. Build out two new test cases (rows 12 & 13) that were previously just
#    blank, bordered placeholder rows.  Start by cloning the formatting of
#    row 11 (borders / alignment / wrap / font) down onto both rows so the
#    new rows look consistent with the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A11:N11").Copy() | Out-Null
$ws.Range("A12:N12").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:N11").Copy() | Out-Null
$ws.Range("A13:N13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Rows.Item(12).RowHeight = 96.6
$ws.Rows.Item(13).RowHeight = 110.4

# Row 12: TC_ERP_Login_003 - invalid username / valid password
$ws.Range("A12").Value2 = "TS_ERP_001"
$ws.Range("B12").Value2 = "Verify the login functionality of the ERP Application Login Page"
$ws.Range("C12").Value2 = "TC_ERP_Login_003"
$ws.Range("D12").Value2 = "Enter invalid user name and one valid password"
$ws.Range("E12").Value2 = "1. Enter invalid username" + [char]10 + "2. Enter valid password" + [char]10 + "3. Click on the Login Button"
$ws.Range("F12").Value2 = "Valid URL Test Data"
$ws.Range("G12").Value2 = "Username: xxxx@erp.com" + [char]10 + "Password: P@asw0rd"
$ws.Range("H12").Value2 = "Still remaining on the login page"
$ws.Range("I12").Value2 = "A popup message box to show an error message ""Invalid Username/Password"""
$ws.Range("J12").Value2 = "Message box got displayed"
$ws.Range("K12").Value2 = "Fail"
$ws.Range("L12").Value2 = "Tester_TLD0001"
$ws.Range("M12").Value2 = 45120
$ws.Range("N12").Value2 = "No comments"

# Row 13: TC_ERP_Login_004 - invalid username / invalid password
$ws.Range("A13").Value2 = "TS_ERP_001"
$ws.Range("B13").Value2 = "Verify the login functionality of the ERP Application Login Page"
$ws.Range("C13").Value2 = "TC_ERP_Login_004"
$ws.Range("D13").Value2 = "Enter invalid user name and one invalid password"
$ws.Range("E13").Value2 = "1. Enter invalid username" + [char]10 + "2. Enter invalid password" + [char]10 + "3. Click on the Login Button"
$ws.Range("F13").Value2 = "Valid URL Test Data"
$ws.Range("G13").Value2 = "Username: xxxx@erp.com" + [char]10 + "Password: xxxxxxxx"
$ws.Range("H13").Value2 = "Still remaining on the login page"
$ws.Range("I13").Value2 = "A popup message box to show an error message ""Invalid Username/Password"""
$ws.Range("J13").Value2 = "Message box got displayed"
$ws.Range("K13").Value2 = "Fail"
$ws.Range("L13").Value2 = "Tester_TLD0001"
$ws.Range("M13").Value2 = 45120
$ws.Range("N13").Value2 = "No comments"



# ---------------------------------------------------------------------------
# 4. Turn the A9:N13 range into a proper Excel Table (ListObject) so the
#    whole test-case grid (header + 4 case rows) is a named table.
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A9:N11"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium1"
$tbl.Resize($ws.Range("A9:N13"))
$tbl.ShowAutoFilter = $false

# ---------------------------------------------------------------------------
# 5. Column width touch-ups so the wider new content (columns B/F/H/J/L/M)
#    fits without truncation - mirrors an "AutoFit columns" pass.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 20.5
$ws.Columns.Item(6).ColumnWidth = 12.67
$ws.Columns.Item(8).ColumnWidth = 13.33
$ws.Columns.Item(10).ColumnWidth = 11.5
$ws.Columns.Item(12).ColumnWidth = 10.67
$ws.Columns.Item(13).ColumnWidth = 14.17

# ---------------------------------------------------------------------------
# 6. Reset view: scroll back to the top-left and select E10 (matches the
#    author's final cursor position after the edit).
# ---------------------------------------------------------------------------
$ws.Range("E10").Select() | Out-Null
